# Commit: "ajout d'une phrase sur le compte rendu"
# Adds a closing sentence to the report, explaining that the user should
# maximise the window when sending a file from a folder whose path is too
# long to fit on the input line. This also relocates the (cursor-tracking)
# "_GoBack" bookmark, and merges the "Description" heading runs that used
# to be split by that bookmark's old position.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the "Description" / " du fonctionnement du programme" runs
#    that are currently split apart by the (old) _GoBack bookmark.
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count

$descPara = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.StartsWith("Description")) {
        $descPara = $p
        break
    }
}

if ($descPara -ne $null) {
    $fullRange = $descPara.Range
    $textLen = $fullRange.Text.Length
    # drop the trailing paragraph mark
    $contentRange = $d.Range($fullRange.Start, $fullRange.Start + $textLen - 1)
    $mergedText = $contentRange.Text

    # Grab the formatting of the first character (red, 20pt heading run)
    # and reapply it across the whole reunified run.
    $srcFormat = $d.Range($fullRange.Start, $fullRange.Start + 1).FormattedText
    $target = $d.Range($fullRange.Start, $fullRange.Start + $textLen - 1)
    $target.FormattedText = $srcFormat

    $newRange = $d.Range($fullRange.Start, $fullRange.Start + 1)
    $newRange.Text = $mergedText
}

# ---------------------------------------------------------------------
# 2) Fill in the final (empty) paragraph with the new closing sentence,
#    and give it the same "spacing after" as the rest of the report.
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$n = $paras.Count
$last = $paras.Item($n)
$last.Format.SpaceAfter = 30

$prev = $paras.Item($n - 1)
$styleSample = $d.Range($prev.Range.Start, $prev.Range.Start + 1)

$apos = [char]0x2019
$eacute = [char]0xE9
$ecirc = [char]0xEA

$part1 = "Si l" + $apos + "utilisateur veut envoyer un fichier qui est dans un dossier dont le chemin ne tient pas sur la ligne d" + $apos + $eacute + "criture, il faut que l" + $apos + "utilisateur ouvre la fen"
$part2 = "$ecirc"
$part3 = "tre en plein " + $eacute + "cran."

$startPos = $last.Range.Start

$r1 = $d.Range($startPos, $startPos)
$r1.FormattedText = $styleSample.FormattedText
$r1 = $d.Range($startPos, $startPos + 1)
$r1.Text = $part1

$pos2 = $startPos + $part1.Length
$r2 = $d.Range($pos2, $pos2)
$r2.FormattedText = $styleSample.FormattedText
$r2 = $d.Range($pos2, $pos2 + 1)
$r2.Text = $part2

# Force a run boundary between part1/part2 (otherwise identically
# formatted adjacent runs get coalesced) using a throwaway bookmark.
$bmWrap = $d.Range($pos2, $pos2 + $part2.Length)
$d.Bookmarks.Add("ZZTempSplit", $bmWrap)

$pos3 = $pos2 + $part2.Length
$r3 = $d.Range($pos3, $pos3)
$r3.FormattedText = $styleSample.FormattedText
$r3 = $d.Range($pos3, $pos3 + 1)
$r3.Text = $part3

$d.Bookmarks.Item("ZZTempSplit").Delete()

# ---------------------------------------------------------------------
# 3) Put the "_GoBack" bookmark back, right where Word would have left
#    it after typing this sentence - collapsed, right after the "ê".
# ---------------------------------------------------------------------
$bmGoBack = $d.Range($pos3, $pos3)
$d.Bookmarks.Add("_GoBack", $bmGoBack)
